# Update countries & provincias Spain
# Applies the 4-Sep-2020 data refresh to the "Pais" sheet:
#  - bumps the "Datos actualizados..." timestamp
#  - refreshes case counters for several countries
#  - the table is kept sorted descending by "Casos totales" (col B), so a
#    few countries swap rows as their totals overtake their neighbours

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 21:32"

# --- Helper: write one data row (country name + 7 numeric columns) ----
function Set-CountryRow {
    param($Row, $Name, $Totales, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)

    $ws.Cells.Item($Row, 1).Value = $Name
    $ws.Cells.Item($Row, 2).Value = $Totales
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Straight data refresh, no re-ranking needed -----------------------
Set-CountryRow 4 "Estados Unidos" 6366921 31677 3585143 2590166 0 554 191612
Set-CountryRow 6 "India"          4020239 87115 3104512  846092 0 1066 69635

# --- Costa Rica overtakes Ghana and Armenia -----------------------------
Set-CountryRow 58 "Costa Rica" 45680 1222 18053 27158 0 9   469
Set-CountryRow 59 "Ghana"      44713 0    43577   856 0 0   280
Set-CountryRow 60 "Armenia"    44461 190  39257  4313 0 4   891

# --- Haiti overtakes Namibia --------------------------------------------
Set-CountryRow 101 "Haiti"   8326 25  5870 2244 0 2 212
Set-CountryRow 102 "Namibia" 8323 241 3611 4625 0 1 87

# --- Cuba overtakes Mozambique, Ruanda and Surinam ----------------------
Set-CountryRow 119 "Cuba"       4266 52 3487  679 0 0 100
Set-CountryRow 120 "Mozambique" 4265 58 2511 1728 0 0 26
Set-CountryRow 121 "Ruanda"     4255 0  2163 2074 0 0 18
Set-CountryRow 122 "Surinam"    4215 0  3318  824 0 0 73

# --- Montserrat swaps with Islas Malvinas (same totals, tie re-break) --
Set-CountryRow 214 "Montserrat"     13 0 12 0 0 0 1
Set-CountryRow 215 "Islas Malvinas" 13 0 13 0 0 0 0
